$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names / emails (shared strings content)
$ws.Range("A1").Value = "claudioneir bossa"
$ws.Range("B1").Value = "claudioneir@gmail.com"
$ws.Range("A2").Value = "Viriosvando bastos"
$ws.Range("B2").Value = "viriosvando@gmail.com"

# Update numeric ids
$ws.Range("D1").Value = 3214569
$ws.Range("D2").Value = 3216544

# Update hyperlinks to match new emails
$ws.Hyperlinks.Item(1).Address = "mailto:claudioneir@gmail.com"
$ws.Hyperlinks.Item(1).TextToDisplay = "claudioneir@gmail.com"
$ws.Hyperlinks.Item(2).Address = "mailto:viriosvando@gmail.com"
$ws.Hyperlinks.Item(2).TextToDisplay = "viriosvando@gmail.com"

# Update selection to D2
$ws.Range("D2").Select()
